$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows B6, B15, B24, B28 currently hold mis-scaled numeric values (e.g. 409275
# instead of 40.93) formatted with a custom "#,##0" style (s="3"). They need
# to become plain text cells (like the rest of column B) holding the corrected
# decimal string. Row B29 holds a bare number (21) with the default style and
# needs to become text "21.0" as well.
#
# To force a text (shared-string) cell instead of Excel re-parsing the string
# back into a number, briefly apply a text number format ("@"), set the
# value, then restore the "Normal" style so the cell ends up with no custom
# style applied (matching the other plain string cells in the column).

function Set-TextValue {
    param($cellAddr, $text)
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "B6" "40.93"
Set-TextValue "B15" "51.3"
Set-TextValue "B24" "41.84"
Set-TextValue "B28" "51.3"
Set-TextValue "B29" "21.0"

# Remove the lingering cell selection stored in the worksheet view (the
# target workbook no longer records an explicit selection).
$ws.Range("A1").Select()
